$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G:K hold numbers written as literal text in the source data
# (e.g. "10", "8", "125.00"). Mark them as Text *before* assigning the
# values so Excel doesn't auto-convert the numeric-looking strings into
# real numbers; columns A:F are already non-numeric-looking and stay text
# on their own, so we leave their formatting untouched.
$ws.Range("G3:K4").NumberFormat = "@"

# Row 3
$ws.Range("A3").Value = " Oct 20 2020"
$ws.Range("B3").Value = " Dubai (DSC)"
$ws.Range("C3").Value = "Kings XI won by 5 wickets (with 6 balls remaining)"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Delhi Capitals"
$ws.Range("F3").Value = "James Neesham "
$ws.Range("G3").Value = "10"
$ws.Range("H3").Value = "8"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "1"
$ws.Range("K3").Value = "125.00"

# Row 4
$ws.Range("A4").Value = " Oct 1 2020"
$ws.Range("B4").Value = " Abu Dhabi"
$ws.Range("C4").Value = "Mumbai won by 48 runs"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "James Neesham "
$ws.Range("G4").Value = "7"
$ws.Range("H4").Value = "7"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "100.00"
